$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resize the table to cover the new rows (A1:C18) ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C18"))

$formula = '=+CONCATENATE("https://www.oanda.com/currency-converter/es/?from=",Tabla1[[#This Row],[TO]],"&to=",Tabla1[[#This Row],[FROM]],"&amount=1")'

# --- Remove the underline formatting from rows 5-7 (B:C) ---
$ws.Range("B5:C7").Font.Underline = -4142

# --- New currency pair rows 8-18 ---
$rows = @(
    @{r=8;  to="EUR"; from="COP"},
    @{r=9;  to="EUR"; from="USD"},
    @{r=10; to="CNY"; from="USD"},
    @{r=11; to="JPY"; from="USD"},
    @{r=12; to="CNY"; from="COP"},
    @{r=13; to="JPY"; from="COP"},
    @{r=14; to="BRL"; from="USD"},
    @{r=15; to="JPY"; from="HNL"},
    @{r=16; to="MXN"; from="HNL"},
    @{r=17; to="HKD"; from="USD"},
    @{r=18; to="HKD"; from="HNL"}
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A${r}").Formula = $formula
    $ws.Range("B${r}").Value = $row.to
    $ws.Range("C${r}").Value = $row.from
}

# Match formatting of the existing table rows (no underline, no fill)
$ws.Range("A8:C18").Font.Underline = -4142

# A13 keeps its original underlined font (matches diff: A13 stays style s="4")
$ws.Range("A13").Font.Underline = 2

# --- Extra formatted (but empty) rows 19-25 below the table ---
$ws.Range("A19:C25").Font.Underline = -4142

# --- Sheet view / selection ---
$ws.Range("D14").Select() | Out-Null

Write-Host "done"
